$d = $word.ActiveDocument
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$r = $d.Range(0,0)
$r.Find.Execute("Need to update line 422 with the search range of the ")
$start = $r.Start
$r.Text = "Need to update line 420 and 411 with the search range of the "

$b1 = $start + ("Need to update line 42").Length
$d.Bookmarks.Add("tmpB1", $d.Range($b1, $b1))
$d.Bookmarks.Item("tmpB1").Delete()

$b2 = $start + ("Need to update line 420 and 411").Length
$d.Bookmarks.Add("_GoBack", $d.Range($b2, $b2))

$b3 = $start + ("Need to update line 420 and 411 with the search range of the ").Length
$d.Bookmarks.Add("tmpB3", $d.Range($b3, $b3))
$d.Bookmarks.Item("tmpB3").Delete()
